$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 61
$ws.Range("H61").Value = 292.6
$ws.Range("I61").Value = 292.6
$ws.Range("K61").Value = 877.8000000000001
$ws.Range("M61").Value = -705.8000000000001

# Row 82
$ws.Range("H82").Value = 6204.125
$ws.Range("I82").Value = 4090.5715
$ws.Range("K82").Value = 12271.7145
$ws.Range("M82").Value = -11865.7145

# Row 85
$ws.Range("H85").Value = 6204.125
$ws.Range("I85").Value = 4090.5715
$ws.Range("K85").Value = 12271.7145
$ws.Range("M85").Value = -10867.7145

# Row 106
$ws.Range("H106").Value = 29335732
$ws.Range("I106").Value = 33848384
$ws.Range("J106").Value = 3499.5
$ws.Range("K106").Value = 33848384
$ws.Range("L106").Value = 3499.5
$ws.Range("M106").Value = -33847753
$ws.Range("N106").Value = -4761.5

# Row 113
$ws.Range("H113").Value = 6500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 6500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6500
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -13008

# Row 132
$ws.Range("H132").Value = 2690.125
$ws.Range("I132").Value = 1962.8
$ws.Range("J132").Value = 3902.3333
$ws.Range("K132").Value = 5888.4
$ws.Range("L132").Value = 11706.9999
$ws.Range("M132").Value = -3358.4
$ws.Range("N132").Value = -16766.9999

# Row 137
$ws.Range("H137").Value = 2840.5652
$ws.Range("I137").Value = 2829.5
$ws.Range("J137").Value = 2880.4
$ws.Range("K137").Value = 8488.5
$ws.Range("L137").Value = 8641.200000000001
$ws.Range("M137").Value = -5938.5
$ws.Range("N137").Value = -13741.2

# Row 138
$ws.Range("H138").Value = 2395.5952
$ws.Range("I138").Value = 1869.6154
$ws.Range("J138").Value = 3250.3125
$ws.Range("K138").Value = 5608.8462
$ws.Range("L138").Value = 9750.9375
$ws.Range("M138").Value = -468.8462
$ws.Range("N138").Value = -20030.9375

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 875.13336
$ws.Range("I2").Value = 901.9286
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 901.9286
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -788.9286
$ws.Range("N2").Value = -726

# Row 61
$ws.Range("H61").Value = 4744.517
$ws.Range("I61").Value = 4151.5186
$ws.Range("J61").Value = 12750
$ws.Range("K61").Value = 4151.5186
$ws.Range("L61").Value = 12750
$ws.Range("M61").Value = -3939.5186
$ws.Range("N61").Value = -13174

# Row 116
$ws.Range("H116").Value = 875.13336
$ws.Range("I116").Value = 901.9286
$ws.Range("J116").Value = 500
$ws.Range("K116").Value = 901.9286
$ws.Range("L116").Value = 500
$ws.Range("M116").Value = 1392.0714
$ws.Range("N116").Value = -5088

# Row 122
$ws.Range("H122").Value = 1582.3889
$ws.Range("I122").Value = 1455
$ws.Range("J122").Value = 3748
$ws.Range("K122").Value = 4365
$ws.Range("L122").Value = 11244
$ws.Range("M122").Value = -1915
$ws.Range("N122").Value = -16144

# Row 136
$ws.Range("H136").Value = 4744.517
$ws.Range("I136").Value = 4151.5186
$ws.Range("J136").Value = 12750
$ws.Range("K136").Value = 12454.5558
$ws.Range("L136").Value = 38250
$ws.Range("M136").Value = -9904.555800000002
$ws.Range("N136").Value = -43350

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 875.13336
$ws.Range("I3").Value = 901.9286
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 901.9286
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = -787.9286
$ws.Range("N3").Value = -728

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 901
$ws.Range("I16").Value = 901
$ws.Range("K16").Value = 901
$ws.Range("M16").Value = -614

# Row 31
$ws.Range("H31").Value = 7547.977
$ws.Range("I31").Value = 4595.273
$ws.Range("J31").Value = 8532.212
$ws.Range("K31").Value = 4595.273
$ws.Range("L31").Value = 8532.212
$ws.Range("M31").Value = -4300.273
$ws.Range("N31").Value = -9122.212

# Row 34
$ws.Range("H34").Value = 7547.977
$ws.Range("I34").Value = 4595.273
$ws.Range("J34").Value = 8532.212
$ws.Range("K34").Value = 4595.273
$ws.Range("L34").Value = 8532.212
$ws.Range("M34").Value = -4393.273
$ws.Range("N34").Value = -8936.212

# Row 107
$ws.Range("H107").Value = 1591.9656
$ws.Range("I107").Value = 1335.48
$ws.Range("K107").Value = 1335.48
$ws.Range("M107").Value = 584.52

# Row 113
$ws.Range("H113").Value = 901
$ws.Range("I113").Value = 901
$ws.Range("K113").Value = 901
$ws.Range("M113").Value = 1269

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 726.8889
$ws.Range("I107").Value = 498
$ws.Range("J107").Value = 910
$ws.Range("K107").Value = 498
$ws.Range("L107").Value = 910
$ws.Range("M107").Value = 1422
$ws.Range("N107").Value = -4750

# Row 113
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 10000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -14340

# Row 122
$ws.Range("H122").Value = 6188.2
$ws.Range("I122").Value = 4441.1816
$ws.Range("J122").Value = 18999.666
$ws.Range("K122").Value = 13323.5448
$ws.Range("L122").Value = 56998.99800000001
$ws.Range("M122").Value = -10873.5448
$ws.Range("N122").Value = -61898.99800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 4994.4165
$ws.Range("J16").Value = 17493.166
$ws.Range("L16").Value = 17493.166
$ws.Range("N16").Value = -17833.166

# Row 68
$ws.Range("H68").Value = 9104.579
$ws.Range("I68").Value = 6454.9
$ws.Range("J68").Value = 12048.667
$ws.Range("K68").Value = 6454.9
$ws.Range("L68").Value = 12048.667
$ws.Range("M68").Value = -5705.9
$ws.Range("N68").Value = -13546.667

# Row 71
$ws.Range("H71").Value = 9104.579
$ws.Range("I71").Value = 6454.9
$ws.Range("J71").Value = 12048.667
$ws.Range("K71").Value = 32274.5
$ws.Range("L71").Value = 60243.335
$ws.Range("M71").Value = -28530.5
$ws.Range("N71").Value = -67731.33499999999

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1883
$ws.Range("I126").Value = 1883
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5649
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3179
$ws.Range("N126").ClearContents()
